# "Add files via upload" — re-upload of class/三年三班.xlsx with updated
# column headers, refreshed survey values, and a different saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: rename the last three columns (shared strings sex/age/grade -> new names)
$ws.Range("E1").Value = "self_piece"
$ws.Range("F1").Value = "win_tour"
$ws.Range("G1").Value = "peer_piece"

# Refresh the survey data in columns F (win_tour) and G (peer_piece) for each student
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0

$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Restore the cursor/selection that was active when the file was last saved
$ws.Range("I16").Select()

# Match the saved application window width (cosmetic Excel chrome state)
$excel.ActiveWindow.Width = 24930
